$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 242
$ws.Range("J9").Value = 259.33334
$ws.Range("L9").Value = 259.33334
$ws.Range("N9").Value = -597.33334
$ws.Range("H15").Value = 772.2258
$ws.Range("I15").Value = 772.2258
$ws.Range("K15").Value = 2316.6774
$ws.Range("M15").Value = -2147.6774
$ws.Range("H86").Value = 15329.583
$ws.Range("I86").Value = 12500
$ws.Range("K86").Value = 12500
$ws.Range("M86").Value = -11377
$ws.Range("H89").Value = 15329.583
$ws.Range("I89").Value = 12500
$ws.Range("K89").Value = 62500
$ws.Range("M89").Value = -56884
$ws.Range("H98").Value = 3005.4
$ws.Range("J98").Value = 4079.1667
$ws.Range("L98").Value = 4079.1667
$ws.Range("N98").Value = -7075.1667
$ws.Range("H100").Value = 3699.75
$ws.Range("I100").Value = 4266.3335
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 4266.3335
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -3725.3335
$ws.Range("N100").Value = -3082
$ws.Range("H106").Value = 3319.8
$ws.Range("I106").Value = 3533
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 3533
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -2902
$ws.Range("N106").Value = -4262
$ws.Range("H112").Value = 2971
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 3156.6667
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 9470.000100000001
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -11686.0001
$ws.Range("H122").Value = 3005.4
$ws.Range("J122").Value = 4079.1667
$ws.Range("L122").Value = 12237.5001
$ws.Range("N122").Value = -17137.5001
$ws.Range("H125").Value = 2349.4285
$ws.Range("I125").Value = 741
$ws.Range("J125").Value = 12000
$ws.Range("K125").Value = 6669
$ws.Range("L125").Value = 108000
$ws.Range("M125").Value = -4209
$ws.Range("N125").Value = -112920

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34997.5
$ws.Range("J44").Value = 34997.5
$ws.Range("L44").Value = 34997.5
$ws.Range("N44").Value = -35973.5
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 30000
$ws.Range("L55").Value = 30000
$ws.Range("N55").Value = -30630
$ws.Range("H132").Value = 2978.3
$ws.Range("I132").Value = 2978.3
$ws.Range("K132").Value = 8934.900000000001
$ws.Range("M132").Value = -6404.900000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 7319.96
$ws.Range("I134").Value = 7713.864
$ws.Range("K134").Value = 23141.592
$ws.Range("M134").Value = -20606.592

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 462.25
$ws.Range("I5").Value = 241.42857
$ws.Range("J5").Value = 2008
$ws.Range("K5").Value = 241.42857
$ws.Range("L5").Value = 2008
$ws.Range("M5").Value = -129.42857
$ws.Range("N5").Value = -2232
$ws.Range("H99").Value = 7383
$ws.Range("I99").Value = 8100
$ws.Range("J99").Value = 6666
$ws.Range("K99").Value = 8100
$ws.Range("L99").Value = 6666
$ws.Range("M99").Value = -6602
$ws.Range("N99").Value = -9662
$ws.Range("H122").Value = 5058.778
$ws.Range("I122").Value = 5058.778
$ws.Range("K122").Value = 15176.334
$ws.Range("M122").Value = -12726.334
$ws.Range("H126").Value = 7383
$ws.Range("I126").Value = 8100
$ws.Range("J126").Value = 6666
$ws.Range("K126").Value = 24300
$ws.Range("L126").Value = 19998
$ws.Range("M126").Value = -21830
$ws.Range("N126").Value = -24938

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 128.34782
$ws.Range("J2").Value = 52.375
$ws.Range("L2").Value = 314.25
$ws.Range("N2").Value = -540.25
$ws.Range("H44").Value = 1004
$ws.Range("I44").Value = 1004
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3012
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -2614
$ws.Range("N44").ClearContents()
$ws.Range("H47").Value = 57.333332
$ws.Range("I47").Value = 57.333332
$ws.Range("K47").Value = 171.999996
$ws.Range("M47").Value = 259.000004
$ws.Range("H55").Value = 350
$ws.Range("I55").Value = 350
$ws.Range("K55").Value = 1050
$ws.Range("M55").Value = -873
$ws.Range("H57").Value = 4593.3335
$ws.Range("I57").Value = 2890
$ws.Range("J57").Value = 8000
$ws.Range("K57").Value = 8670
$ws.Range("L57").Value = 24000
$ws.Range("M57").Value = -8111
$ws.Range("N57").Value = -25118
$ws.Range("H59").Value = 268.33334
$ws.Range("J59").Value = 400
$ws.Range("L59").Value = 1200
$ws.Range("N59").Value = -2280
$ws.Range("H60").Value = 574.8
$ws.Range("I60").Value = 218.5
$ws.Range("J60").Value = 2000
$ws.Range("K60").Value = 655.5
$ws.Range("L60").Value = 6000
$ws.Range("M60").Value = -404.5
$ws.Range("N60").Value = -6502
$ws.Range("H138").Value = 638.3333
$ws.Range("I138").Value = 638.3333
$ws.Range("K138").Value = 1914.9999
$ws.Range("M138").Value = 3225.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3063.6
$ws.Range("I102").Value = 3063.6
$ws.Range("K102").Value = 3063.6
$ws.Range("M102").Value = -1441.6
$ws.Range("H122").Value = 49786.734
$ws.Range("I122").Value = 54708.453
$ws.Range("K122").Value = 164125.359
$ws.Range("M122").Value = -161675.359

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5698.7144
$ws.Range("I40").Value = 5832
$ws.Range("J40").Value = 4899
$ws.Range("K40").Value = 5832
$ws.Range("L40").Value = 4899
$ws.Range("M40").Value = -5696
$ws.Range("N40").Value = -5171
$ws.Range("H132").Value = 12119.77
$ws.Range("I132").Value = 16209.111
$ws.Range("J132").Value = 2918.75
$ws.Range("K132").Value = 48627.333
$ws.Range("L132").Value = 8756.25
$ws.Range("M132").Value = -46097.333
$ws.Range("N132").Value = -13816.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H136").Value = 1989.4117
$ws.Range("I136").Value = 1863.75
$ws.Range("K136").Value = 5591.25
$ws.Range("M136").Value = -3041.25
